$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.943.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.289.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.118"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.436"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.840.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.12%  "
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000180"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.731.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.281.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.14%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.534"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.174"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0977"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.07%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.64%  "
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.55%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.34%  "
$ws.Range("E34").Value = "  +4.82%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.40%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.84%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.822.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0726"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0319"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.84%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.740"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.33%  "
$ws.Range("E45").Value = "  +4.33%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.28%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.105"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.321.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.811"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "280.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.37%  "
